{"js": "// Replace each old equation text with its new equation text.\n// Mapping taken 1:1 from the diff: each left-hand side text is unique in\n// the document, so a matchCase, non-wildcard body.search() + insertText\n// 'Replace' on every hit is unambiguous.\nconst replacements = [\n  [\"868\u00d77=\", \"121\u00d76=\"],\n  [\"711\u00d75=\", \"123\u00d76=\"],\n  [\"205\u00d78=\", \"857\u00d79=\"],\n  [\"646\u00d79=\", \"384\u00d75=\"],\n  [\"678\u00d73=\", \"273\u00d75=\"],\n  [\"493\u00d76=\", \"877\u00d75=\"],\n  [\"723\u00d78=\", \"226\u00d75=\"],\n  [\"347\u00d74=\", \"702\u00d74=\"],\n  [\"596\u00d77=\", \"252\u00d72=\"],\n  [\"171\u00d78=\", \"680\u00d73=\"],\n  [\"697\u00d75=\", \"473\u00d78=\"],\n  [\"204\u00d79=\", \"239\u00d79=\"],\n  [\"972\u00d78=\", \"129\u00d72=\"],\n  [\"599\u00d76=\", \"726\u00d79=\"],\n  [\"235\u00d74=\", \"787\u00d72=\"],\n  [\"466\u00d76=\", \"826\u00d77=\"],\n  [\"602\u00d77=\", \"797\u00d76=\"],\n  [\"828\u00d78=\", \"274\u00d72=\"],\n  [\"609\u00d73=\", \"803\u00d75=\"],\n  [\"574\u00d72=\", \"953\u00d74=\"],\n  [\"181\u00d76=\", \"291\u00d73=\"],\n  [\"601\u00d76=\", \"803\u00d73=\"],\n  [\"911\u00d74=\", \"126\u00d79=\"],\n  [\"475\u00d74=\", \"176\u00d79=\"],\n  [\"581\u00d75=\", \"319\u00d78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const result of results.items) {\n    result.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each old equation text with its new equation text.\n# Each left-hand side text is unique in the document, so a single\n# Find/Replace (wdReplaceAll) pass per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"868\u00d77=\", \"121\u00d76=\"),\n    @(\"711\u00d75=\", \"123\u00d76=\"),\n    @(\"205\u00d78=\", \"857\u00d79=\"),\n    @(\"646\u00d79=\", \"384\u00d75=\"),\n    @(\"678\u00d73=\", \"273\u00d75=\"),\n    @(\"493\u00d76=\", \"877\u00d75=\"),\n    @(\"723\u00d78=\", \"226\u00d75=\"),\n    @(\"347\u00d74=\", \"702\u00d74=\"),\n    @(\"596\u00d77=\", \"252\u00d72=\"),\n    @(\"171\u00d78=\", \"680\u00d73=\"),\n    @(\"697\u00d75=\", \"473\u00d78=\"),\n    @(\"204\u00d79=\", \"239\u00d79=\"),\n    @(\"972\u00d78=\", \"129\u00d72=\"),\n    @(\"599\u00d76=\", \"726\u00d79=\"),\n    @(\"235\u00d74=\", \"787\u00d72=\"),\n    @(\"466\u00d76=\", \"826\u00d77=\"),\n    @(\"602\u00d77=\", \"797\u00d76=\"),\n    @(\"828\u00d78=\", \"274\u00d72=\"),\n    @(\"609\u00d73=\", \"803\u00d75=\"),\n    @(\"574\u00d72=\", \"953\u00d74=\"),\n    @(\"181\u00d76=\", \"291\u00d73=\"),\n    @(\"601\u00d76=\", \"803\u00d73=\"),\n    @(\"911\u00d74=\", \"126\u00d79=\"),\n    @(\"475\u00d74=\", \"176\u00d79=\"),\n    @(\"581\u00d75=\", \"319\u00d78=\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n    #         MatchSoundsLike, MatchAllWordForms, Forward, Wrap,\n    #         Format, ReplaceWith, Replace)\n    # Wrap=1 (wdFindContinue), Replace=2 (wdReplaceAll)\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
